$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (OKLO)
$ws.Range("D2").Value = 104.91
$ws.Range("E2").Value = 55.2
$ws.Range("F2").Value = 14.81
$ws.Range("K2").Value = 63.1
$ws.Range("N2").Value = 49.16024380385575

# Row 3 (SMR)
$ws.Range("D3").Value = 21.93
$ws.Range("E3").Value = 48.4
$ws.Range("F3").Value = 9.65
$ws.Range("H3").Value = 56
$ws.Range("K3").Value = 54.9
$ws.Range("N3").Value = 49.16024380385575
